$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text casing for Unit SN column
$ws.Range("A1").Value = "Unit SN"

# Update the active selection to match the new cursor position
$ws.Range("G4").Select()
